$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Taxa" column (B) values to 0 for rows 2-10
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0

# Update the "Data de Salvamento" column (C) for all data rows (2-15)
# with the new timestamp as text
for ($row = 2; $row -le 15; $row++) {
    $ws.Range("C$row").Value = "2025-04-04 14:29:21"
}
